$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 94
$ws.Cells.Item(2, 3).Value = "flower/flower106.png"
$ws.Cells.Item(2, 4).Value = "fliegen"
$ws.Cells.Item(2, 5).Value = "flower"
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = "flower/flower089.png"
$ws.Cells.Item(3, 4).Value = "bleiben"
$ws.Cells.Item(3, 5).Value = "flower"
$ws.Cells.Item(4, 2).Value = 73
$ws.Cells.Item(4, 3).Value = "flower/flower083.png"
$ws.Cells.Item(4, 4).Value = "schätzen"
$ws.Cells.Item(4, 5).Value = "flower"
$ws.Cells.Item(5, 2).Value = 121
$ws.Cells.Item(5, 3).Value = "flower/flower100.png"
$ws.Cells.Item(5, 4).Value = "loben"
$ws.Cells.Item(5, 5).Value = "flower"
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = "flower/flower095.png"
$ws.Cells.Item(6, 4).Value = "formen"
$ws.Cells.Item(6, 5).Value = "flower"
$ws.Cells.Item(7, 2).Value = 124
$ws.Cells.Item(7, 3).Value = "dog/dog068.png"
$ws.Cells.Item(7, 4).Value = "dauern"
$ws.Cells.Item(7, 5).Value = "dog"
$ws.Cells.Item(8, 2).Value = 8
$ws.Cells.Item(8, 3).Value = "dog/dog064.png"
$ws.Cells.Item(8, 4).Value = "regnen"
$ws.Cells.Item(8, 5).Value = "dog"
$ws.Cells.Item(9, 2).Value = 42
$ws.Cells.Item(9, 3).Value = "dog/dog091.png"
$ws.Cells.Item(9, 4).Value = "fliehen"
$ws.Cells.Item(9, 5).Value = "dog"
$ws.Cells.Item(10, 2).Value = 79
$ws.Cells.Item(10, 3).Value = "dog/dog109.png"
$ws.Cells.Item(10, 4).Value = "töten"
$ws.Cells.Item(10, 5).Value = "dog"
$ws.Cells.Item(11, 2).Value = 20
$ws.Cells.Item(11, 3).Value = "flower/flower091.png"
$ws.Cells.Item(11, 4).Value = "krachen"
$ws.Cells.Item(11, 5).Value = "flower"
$ws.Cells.Item(12, 2).Value = 51
$ws.Cells.Item(12, 3).Value = "flower/flower070.png"
$ws.Cells.Item(12, 4).Value = "antun"
$ws.Cells.Item(12, 5).Value = "flower"
$ws.Cells.Item(13, 2).Value = 117
$ws.Cells.Item(13, 3).Value = "dog/dog114.png"
$ws.Cells.Item(13, 4).Value = "spielen"
$ws.Cells.Item(13, 5).Value = "dog"
$ws.Cells.Item(14, 2).Value = 103
$ws.Cells.Item(14, 3).Value = "flower/flower110.png"
$ws.Cells.Item(14, 4).Value = "posten"
$ws.Cells.Item(14, 5).Value = "flower"
$ws.Cells.Item(15, 2).Value = 102
$ws.Cells.Item(15, 3).Value = "flower/flower082.png"
$ws.Cells.Item(15, 4).Value = "stärken"
$ws.Cells.Item(15, 5).Value = "flower"
$ws.Cells.Item(16, 2).Value = 101
$ws.Cells.Item(16, 3).Value = "dog/dog066.png"
$ws.Cells.Item(16, 4).Value = "runden"
$ws.Cells.Item(16, 5).Value = "dog"
$ws.Cells.Item(17, 2).Value = 65
$ws.Cells.Item(17, 3).Value = "flower/flower073.png"
$ws.Cells.Item(17, 4).Value = "tauschen"
$ws.Cells.Item(17, 5).Value = "flower"
$ws.Cells.Item(18, 2).Value = 82
$ws.Cells.Item(18, 3).Value = "dog/dog074.png"
$ws.Cells.Item(18, 4).Value = "stechen"
$ws.Cells.Item(18, 5).Value = "dog"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = "flower/flower067.png"
$ws.Cells.Item(19, 4).Value = "wiegen"
$ws.Cells.Item(19, 5).Value = "flower"
$ws.Cells.Item(20, 2).Value = 44
$ws.Cells.Item(20, 3).Value = "flower/flower084.png"
$ws.Cells.Item(20, 4).Value = "gelten"
$ws.Cells.Item(20, 5).Value = "flower"
$ws.Cells.Item(21, 2).Value = 9
$ws.Cells.Item(21, 3).Value = "dog/dog088.png"
$ws.Cells.Item(21, 4).Value = "nehmen"
$ws.Cells.Item(21, 5).Value = "dog"
$ws.Cells.Item(22, 2).Value = 112
$ws.Cells.Item(22, 3).Value = "dog/dog071.png"
$ws.Cells.Item(22, 4).Value = "gründen"
$ws.Cells.Item(22, 5).Value = "dog"
$ws.Cells.Item(23, 2).Value = 40
$ws.Cells.Item(23, 3).Value = "dog/dog085.png"
$ws.Cells.Item(23, 4).Value = "währen"
$ws.Cells.Item(23, 5).Value = "dog"
$ws.Cells.Item(24, 2).Value = 81
$ws.Cells.Item(24, 3).Value = "dog/dog098.png"
$ws.Cells.Item(24, 4).Value = "rücken"
$ws.Cells.Item(24, 5).Value = "dog"
$ws.Cells.Item(25, 2).Value = 114
$ws.Cells.Item(25, 3).Value = "flower/flower117.png"
$ws.Cells.Item(25, 4).Value = "lehnen"
$ws.Cells.Item(25, 5).Value = "flower"
$ws.Cells.Item(26, 2).Value = 89
$ws.Cells.Item(26, 3).Value = "flower/flower080.png"
$ws.Cells.Item(26, 4).Value = "strahlen"
$ws.Cells.Item(26, 5).Value = "flower"
$ws.Cells.Item(27, 2).Value = 45
$ws.Cells.Item(27, 3).Value = "dog/dog073.png"
$ws.Cells.Item(27, 4).Value = "füllen"
$ws.Cells.Item(27, 5).Value = "dog"
$ws.Cells.Item(28, 2).Value = 43
$ws.Cells.Item(28, 3).Value = "dog/dog104.png"
$ws.Cells.Item(28, 4).Value = "fühlen"
$ws.Cells.Item(28, 5).Value = "dog"
$ws.Cells.Item(29, 2).Value = 50
$ws.Cells.Item(29, 3).Value = "flower/flower086.png"
$ws.Cells.Item(29, 4).Value = "backen"
$ws.Cells.Item(29, 5).Value = "flower"
$ws.Cells.Item(30, 2).Value = 18
$ws.Cells.Item(30, 3).Value = "dog/dog084.png"
$ws.Cells.Item(30, 4).Value = "sieben"
$ws.Cells.Item(30, 5).Value = "dog"
$ws.Cells.Item(31, 2).Value = 54
$ws.Cells.Item(31, 3).Value = "dog/dog093.png"
$ws.Cells.Item(31, 4).Value = "füttern"
$ws.Cells.Item(31, 5).Value = "dog"
$ws.Cells.Item(32, 2).Value = 60
$ws.Cells.Item(32, 3).Value = "dog/dog090.png"
$ws.Cells.Item(32, 4).Value = "enden"
$ws.Cells.Item(32, 5).Value = "dog"
$ws.Cells.Item(33, 2).Value = 12
$ws.Cells.Item(33, 3).Value = "flower/flower092.png"
$ws.Cells.Item(33, 4).Value = "raten"
$ws.Cells.Item(33, 5).Value = "flower"
